$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15").Value = "some new descriptive text"
$ws.Range("C15").Formula = '=CONCATENATE("Here is some text to append: ", E15)'
$ws.Columns("C").ColumnWidth = 50.8
$ws.Range("C20").Select()
